$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "2025/12/04 03:00"
$ws.Range("B35").Value = "-"
$ws.Range("C35").Value = "-"
$ws.Range("D35").Value = "-"
$ws.Range("E35").Value = "-"
$ws.Range("F35").Value = "-"
$ws.Range("G35").Value = "-"
